$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 already exists - update Product ID, Price and Visibility; Product Name/Packing Unit ID also refreshed
$ws.Range("A2").Value = 5151
$ws.Range("B2").Value = "ريد بل - 250 مل"
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 1065
$ws.Range("E2").Value = "YES"

$ws.Range("A3").Value = 5152
$ws.Range("B3").Value = "ريد بل فرى شوجر - 250 مل"
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 1065
$ws.Range("E3").Value = "YES"

$ws.Range("A4").Value = 5153
$ws.Range("B4").Value = "ريد بل ابيض بجوز الهند و التوت - 250 مل"
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = 1065
$ws.Range("E4").Value = "YES"

$ws.Range("A5").Value = 11509
$ws.Range("B5").Value = "ريد بل توت ازرق - 250 مل"
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = 1065
$ws.Range("E5").Value = "YES"

$ws.Range("A6").Value = 13928
$ws.Range("B6").Value = "ريد بول 12 كانز - 250 مل"
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 532
$ws.Range("E6").Value = "YES"

$ws.Range("A7").Value = 13928
$ws.Range("B7").Value = "ريد بول 12 كانز - 250 مل"
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = 1065
$ws.Range("E7").Value = "YES"

$ws.Range("A8").Value = 19977
$ws.Range("B8").Value = "ريدبل كريز و توت بري  - 250 مل"
$ws.Range("C8").Value = 2
$ws.Range("D8").Value = 1050
$ws.Range("E8").Value = "YES"

$ws.Range("A9").Value = 7630
$ws.Range("B9").Value = "فيورى جولد - 400 مل"
$ws.Range("C9").Value = 2
$ws.Range("D9").Value = 205
$ws.Range("E9").Value = "YES"

